$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Row 4 (Sr.No 3) already exists: TC_03, Prathamesh, Lad, Chrome, Homepage, Smoke, Yes
# Add 4 more rows (5-8) following the same pattern, taking all test cases
# (not just the ones marked "Yes" for Regression).

$data = @(
    @(4, "TC_03", "Prathamesh", "Lad", "Chrome", "Homepage", "Smoke", "Yes"),
    @(5, "TC_03", "Prathamesh", "Lad", "Chrome", "Homepage", "Yes",   "Yes"),
    @(6, "TC_03", "Prathamesh", "Lad", "Chrome", "Homepage", "Smoke", "Yes"),
    @(7, "TC_03", "Prathamesh", "Lad", "Chrome", "Homepage", "Smoke", "Yes")
)

$r = 5
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}

$ws.Range("H13").Select()
